# Append a new data row (row 93) to Sheet1, mirroring the existing
# "Date, Original (CNY/mt), VAT Included (USD/mt), VAT Excluded (USD/mt), USD/CNY"
# layout used by the preceding rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A93").Value = "2024-11-03 00:00:00"
$ws.Range("B93").Value = 73850
$ws.Range("C93").Value = 10340.96
$ws.Range("D93").Value = 9151.299999999999
$ws.Range("E93").Value = 7.1227
